# processing metals.xlsx - "CrewAI Robust Backend Ready!" edit
#
# 1. Lower-case the header row text (A1:G1)
# 2. Add header-cell comments describing each column's data type
# 3. Shift the numeric columns D/E/F one position to the left
#    (new D = old E, new E = old F) and populate F with the new
#    "climate change (kg CO2 eq)" values
# 4. A legacyDrawing relationship is produced automatically by Excel
#    once cell comments exist on the sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1. Header text -> lower case
# ---------------------------------------------------------------
$ws.Cells.Item(1, 1).Value2 = "industry"
$ws.Cells.Item(1, 2).Value2 = "unit"
$ws.Cells.Item(1, 3).Value2 = "process"
$ws.Cells.Item(1, 4).Value2 = "carbon (kg CO2 eq)"
$ws.Cells.Item(1, 5).Value2 = "ced (MJ)"
$ws.Cells.Item(1, 6).Value2 = "climate change (kg CO2 eq)"
$ws.Cells.Item(1, 7).Value2 = "region"

# ---------------------------------------------------------------
# 2. Header comments (data type annotations)
# ---------------------------------------------------------------
$ws.Cells.Item(1, 1).AddComment("Data type: Categorical (text)")
$ws.Cells.Item(1, 2).AddComment("Data type: Various (e.g. kg, kWh)")
$ws.Cells.Item(1, 3).AddComment("Data type: Categorical (text)")
$ws.Cells.Item(1, 4).AddComment("Data type: Carbon footprint")
$ws.Cells.Item(1, 5).AddComment("Data type: Cumulative energy demand")
$ws.Cells.Item(1, 6).AddComment("Data type: Climate change impact")
$ws.Cells.Item(1, 7).AddComment("Data type: Categorical (text)")

# ---------------------------------------------------------------
# 3. Shift D/E/F columns: new D = old E, new E = old F,
#    new F = freshly computed "climate change" values
# ---------------------------------------------------------------
$newF = @(0.000013880092, 0.000004182398, 0.000000015128166, 0.0000007564083, 0.000002912172, 0.000002912172, 0.000002912172, 0.00006111540900000001, 0.000035092739, 0.000080981879, 0.000040596837, 0.00039775589, 0.000018565961, 0.0000003710246, 0.00000018420623, 0.00011084392, 0.00011029482, 0.000035961236, 0.0058553232, 0.0000063010355, 0.0000051777789, 0.000060454147, 0.000016444317, 0.00004485123, 0.000070266549, 0.0000027684544, 0.000000022692249, 0.0000013615349, 0.0000018153799, 0.000021801084, 0.0000060890868, 0.0000018153799, 0.000092690274, 0.00016296817, 0.000017942005, 0.000038879387, 0.00005831908, 0.00006921136, 0.00008819720800000001, 0.0000018153799, 0.000027701117, 0.000023053698, 0.000012318155, 0.000011213753, 0.000022427506, 0.000023524298, 0.000074752051, 0.00023920656, 0.00034386321, 0.00000074884422, 0.000001495041, 0.0000029915948, 0.0000070270331, 0.0000018153799, 0.0000068131846)

for ($row = 2; $row -le 56; $row++) {
    $oldE = $ws.Cells.Item($row, 5).Value2
    $oldF = $ws.Cells.Item($row, 6).Value2

    $ws.Cells.Item($row, 4).Value2 = $oldE
    $ws.Cells.Item($row, 5).Value2 = $oldF
    $ws.Cells.Item($row, 6).Value2 = $newF[$row - 2]
}
